# Update TPM-derived NATMI ligand/receptor metrics with newly computed TPM values.
# All source columns (A-F) are unchanged; only derived metric columns G-T shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "0.3410236666666666"
$ws.Range("H2").Value = "1.023071"
$ws.Range("I2").Value = "0.01850325494520333"
$ws.Range("J2").Value = "0.01850325494520333"
$ws.Range("K2").Value = "2"
$ws.Range("L2").Value = "0.6666666666666666"
$ws.Range("M2").Value = "0.022105"
$ws.Range("N2").Value = "0.066315"
$ws.Range("O2").Value = "0.0007557226718989593"
$ws.Range("P2").Value = "0.0007557226718989592"
$ws.Range("Q2").Value = "0.007538328151666666"
$ws.Range("R2").Value = "0.06784495336499999"
$ws.Range("S2").Value = [double]"1.398332926601669E-05"
$ws.Range("T2").Value = [double]"1.398332926601669E-05"

$ws.Range("G3").Value = "0.3410236666666666"
$ws.Range("H3").Value = "1.023071"
$ws.Range("I3").Value = "0.01850325494520333"
$ws.Range("J3").Value = "0.01850325494520333"
$ws.Range("O3").Value = "0.01371106452749117"
$ws.Range("P3").Value = "0.01371106452749117"
$ws.Range("Q3").Value = "0.1367677688657778"
$ws.Range("R3").Value = "1.230909919792"
$ws.Range("S3").Value = "0.0002536993225223028"
$ws.Range("T3").Value = "0.0002536993225223028"

$ws.Range("G4").Value = "0.3410236666666666"
$ws.Range("H4").Value = "1.023071"
$ws.Range("I4").Value = "0.01850325494520333"
$ws.Range("J4").Value = "0.01850325494520333"
$ws.Range("M4").Value = "28.82699233333333"
$ws.Range("N4").Value = "86.480977"
$ws.Range("O4").Value = "0.9855332128006099"
$ws.Range("P4").Value = "0.9855332128006098"
$ws.Range("Q4").Value = "9.830686624485221"
$ws.Range("R4").Value = "88.47617962036698"
$ws.Range("S4").Value = "0.01823557229341501"
$ws.Range("T4").Value = "0.01823557229341501"

$ws.Range("I5").Value = "0.2085050756621187"
$ws.Range("J5").Value = "0.2085050756621187"
$ws.Range("K5").Value = "2"
$ws.Range("L5").Value = "0.6666666666666666"
$ws.Range("M5").Value = "0.022105"
$ws.Range("N5").Value = "0.066315"
$ws.Range("O5").Value = "0.0007557226718989593"
$ws.Range("P5").Value = "0.0007557226718989592"
$ws.Range("Q5").Value = "0.08494611819833334"
$ws.Range("R5").Value = "0.764515063785"
$ws.Range("S5").Value = "0.000157572012883871"
$ws.Range("T5").Value = "0.000157572012883871"

$ws.Range("I6").Value = "0.2085050756621187"
$ws.Range("J6").Value = "0.2085050756621187"
$ws.Range("O6").Value = "0.01371106452749117"
$ws.Range("P6").Value = "0.01371106452749117"
$ws.Range("S6").Value = "0.002858826546712737"
$ws.Range("T6").Value = "0.002858826546712737"

$ws.Range("I7").Value = "0.2085050756621187"
$ws.Range("J7").Value = "0.2085050756621187"
$ws.Range("M7").Value = "28.82699233333333"
$ws.Range("N7").Value = "86.480977"
$ws.Range("O7").Value = "0.9855332128006099"
$ws.Range("P7").Value = "0.9855332128006098"
$ws.Range("Q7").Value = "110.7777017891781"
$ws.Range("R7").Value = "996.999316102603"
$ws.Range("S7").Value = "0.2054886771025221"
$ws.Range("T7").Value = "0.2054886771025221"

$ws.Range("G8").Value = "0.2092423333333333"
$ws.Range("H8").Value = "0.627727"
$ws.Range("I8").Value = "0.0113530661283407"
$ws.Range("J8").Value = "0.0113530661283407"
$ws.Range("K8").Value = "2"
$ws.Range("L8").Value = "0.6666666666666666"
$ws.Range("M8").Value = "0.022105"
$ws.Range("N8").Value = "0.066315"
$ws.Range("O8").Value = "0.0007557226718989593"
$ws.Range("P8").Value = "0.0007557226718989592"
$ws.Range("Q8").Value = "0.004625301778333334"
$ws.Range("R8").Value = "0.041627716005"
$ws.Range("S8").Value = [double]"8.579769468755208E-06"
$ws.Range("T8").Value = [double]"8.579769468755206E-06"

$ws.Range("G9").Value = "0.2092423333333333"
$ws.Range("H9").Value = "0.627727"
$ws.Range("I9").Value = "0.0113530661283407"
$ws.Range("J9").Value = "0.0113530661283407"
$ws.Range("O9").Value = "0.01371106452749117"
$ws.Range("P9").Value = "0.01371106452749117"
$ws.Range("Q9").Value = "0.08391677727822222"
$ws.Range("R9").Value = "0.755250995504"
$ws.Range("S9").Value = "0.0001556626222705536"
$ws.Range("T9").Value = "0.0001556626222705536"

$ws.Range("G10").Value = "0.2092423333333333"
$ws.Range("H10").Value = "0.627727"
$ws.Range("I10").Value = "0.0113530661283407"
$ws.Range("J10").Value = "0.0113530661283407"
$ws.Range("M10").Value = "28.82699233333333"
$ws.Range("N10").Value = "86.480977"
$ws.Range("O10").Value = "0.9855332128006099"
$ws.Range("P10").Value = "0.9855332128006098"
$ws.Range("Q10").Value = "6.031827138808778"
$ws.Range("R10").Value = "54.286444249279"
$ws.Range("S10").Value = "0.01118882373660139"
$ws.Range("T10").Value = "0.01118882373660139"

$ws.Range("G11").Value = "14.03735666666667"
$ws.Range("H11").Value = "42.11207"
$ws.Range("I11").Value = "0.7616386032643372"
$ws.Range("J11").Value = "0.7616386032643372"
$ws.Range("K11").Value = "2"
$ws.Range("L11").Value = "0.6666666666666666"
$ws.Range("M11").Value = "0.022105"
$ws.Range("N11").Value = "0.066315"
$ws.Range("O11").Value = "0.0007557226718989593"
$ws.Range("P11").Value = "0.0007557226718989592"
$ws.Range("Q11").Value = "0.3102957691166667"
$ws.Range("R11").Value = "2.79266192205"
$ws.Range("S11").Value = "0.0005755875602803163"
$ws.Range("T11").Value = "0.0005755875602803163"

$ws.Range("G12").Value = "14.03735666666667"
$ws.Range("H12").Value = "42.11207"
$ws.Range("I12").Value = "0.7616386032643372"
$ws.Range("J12").Value = "0.7616386032643372"
$ws.Range("O12").Value = "0.01371106452749117"
$ws.Range("P12").Value = "0.01371106452749117"
$ws.Range("Q12").Value = "5.629691249404445"
$ws.Range("R12").Value = "50.66722124464"
$ws.Range("S12").Value = "0.01044287603598557"
$ws.Range("T12").Value = "0.01044287603598557"

$ws.Range("G13").Value = "14.03735666666667"
$ws.Range("H13").Value = "42.11207"
$ws.Range("I13").Value = "0.7616386032643372"
$ws.Range("J13").Value = "0.7616386032643372"
$ws.Range("M13").Value = "28.82699233333333"
$ws.Range("N13").Value = "86.480977"
$ws.Range("O13").Value = "0.9855332128006099"
$ws.Range("P13").Value = "0.9855332128006098"
$ws.Range("Q13").Value = "404.6547730102656"
$ws.Range("R13").Value = "3641.89295709239"
$ws.Range("S13").Value = "0.7506201396680713"
$ws.Range("T13").Value = "0.7506201396680713"

